$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Bad Photos"

# Update header cell values
$ws.Range("B1").Value = "Local file"
$ws.Range("C1").Value = "Width"
$ws.Range("D1").Value = "Height"
